$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.498.74'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '1.920.73'
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("E4").Value = '  +0.88%  '
$ws.Range("D5").Value = '''326.33'
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").Value = '''1.008'
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("D7").Value = '''0.4832'
$ws.Range("E7").Value = '  +2.54%  '
$ws.Range("D8").Value = '''0.4086'
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("D9").Value = '''0.08231'
$ws.Range("E9").Value = '  +2.59%  '
$ws.Range("E10").Value = '  +2.92%  '
$ws.Range("D11").Value = '''23.54'
$ws.Range("E11").Value = '  +2.70%  '
$ws.Range("D12").Value = '1.938.62'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '''6.051'
$ws.Range("E13").Value = '  +1.81%  '
$ws.Range("D14").Value = '''7.241'
$ws.Range("E14").Value = '  +2.73%  '
$ws.Range("D15").Value = '''91.40'
$ws.Range("E15").Value = '  +2.22%  '
$ws.Range("D16").Value = '''0.06810'
$ws.Range("E16").Value = '  +2.64%  '
$ws.Range("D17").Value = '''1.009'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").Value = '''0.00001042'
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("D19").Value = '''17.83'
$ws.Range("E19").Value = '  +2.08%  '
$ws.Range("D20").Value = '''1.008'
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("D21").Value = '29.523.46'
$ws.Range("E21").Value = '  +0.85%  '
$ws.Range("D22").Value = '''5.650'
$ws.Range("E22").Value = '  +2.89%  '
$ws.Range("E23").Value = '  +0.60%  '
$ws.Range("D24").Value = '''2.196'
$ws.Range("E24").Value = '  +0.98%  '
$ws.Range("D25").Value = '2.165.65'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("D26").Value = '''6.687'
$ws.Range("E26").Value = '  +11.77%  '
$ws.Range("D27").Value = '''156.87'
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").Value = '''20.10'
$ws.Range("E28").Value = '  +2.13%  '
$ws.Range("D29").Value = '''2.121'
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("D30").Value = '''120.58'
$ws.Range("E30").Value = '  +2.92%  '
$ws.Range("D31").Value = '''1.028'
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("D32").Value = '''0.09582'
$ws.Range("E32").Value = '  +1.64%  '
$ws.Range("D33").Value = '''5.543'
$ws.Range("E33").Value = '  +3.46%  '
$ws.Range("D34").Value = '''3.564'
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("D35").Value = '''1.385'
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("D36").Value = '''0.02287'
$ws.Range("E36").Value = '  +1.99%  '
$ws.Range("D37").Value = '''0.06142'
$ws.Range("E37").Value = '  +1.57%  '
$ws.Range("D38").Value = '''1.180'
$ws.Range("E38").Value = '  +0.72%  '
$ws.Range("D39").Value = '''0.5996'
$ws.Range("E39").Value = '  +2.94%  '
$ws.Range("D40").Value = '''8.056'
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("D41").Value = '''10.83'
$ws.Range("E41").Value = '  +7.71%  '
$ws.Range("D42").Value = '''0.1858'
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''2.408'
$ws.Range("E43").Value = '  +1.66%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '''1.281'
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("D45").Value = '''0.07606'
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").Value = '''12.44'
$ws.Range("E46").Value = '  +1.79%  '
$ws.Range("D47").Value = '''0.5589'
$ws.Range("E47").Value = '  +2.05%  '
$ws.Range("D48").Value = '''1.962'
$ws.Range("E48").Value = '  +2.85%  '
$ws.Range("D49").Value = '''118.03'
$ws.Range("E49").Value = '  +4.24%  '
$ws.Range("D50").Value = '''2.430'
$ws.Range("E50").Value = '  +4.57%  '
$ws.Range("D51").Value = '''72.92'
$ws.Range("E51").Value = '  +2.53%  '
